# Add a new "2021" column (O) to the poverty-rate table, mirroring the
# formatting already used for column N (2020), but in the slightly larger
# 10pt "Times New Roman CE" font that the author used for the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: thin bottom-border spacer cell, no value, just style ---
$ws.Range("N3").Copy()
$ws.Range("O3").PasteSpecial(-4122)   # xlPasteFormats

# --- Row 4: bold header "2021" (mirrors N4, font bumped 9pt -> 10pt) ---
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("O4").Value = 2021
$ws.Range("O4").Font.Size = 10

# --- Row 5: first data row (mirrors N5, font bumped 9pt -> 10pt) ---
$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$ws.Range("O5").Value = 6.0337796775071091
$ws.Range("O5").Font.Size = 10

# --- Rows 6-15: remaining data rows (mirror N6-N15, font bumped to 10pt) ---
$values = @{
    6  = 7.3075058743442511
    7  = 5.2767607763499562
    8  = 10.064200140319592
    9  = 7.5445007460298559
    10 = 7.9562092224762884
    11 = 8.1696953402867685
    12 = 2.0701729813092102
    13 = 2.6482523478927704
    14 = 3.9561647100749857
    15 = 9.4645167179465837
}

foreach ($r in 6..15) {
    $ws.Range("N$r").Copy()
    $ws.Range("O$r").PasteSpecial(-4122)
    $ws.Range("O$r").Value = $values[$r]
    $ws.Range("O$r").Font.Size = 10
}

# --- Row 16: bottom total row, thick bottom border (mirrors N16) ---
$ws.Range("N16").Copy()
$ws.Range("O16").PasteSpecial(-4122)
$ws.Range("O16").Value = 3.1019579996103404
$ws.Range("O16").Font.Size = 10

$excel.CutCopyMode = $false

# --- Row heights: rows without an explicit height pick up 12.75pt once the
#     new 10pt column is added; the bottom total row grows from 12.75 to 13.5 ---
foreach ($r in 5,6,7,8,9,10,11,13,14,15) {
    $ws.Rows($r).RowHeight = 12.75
}
$ws.Rows(16).RowHeight = 13.5
